# edit.ps1
# Applies the commit: "Adiciona relatorio 5 e modifica simulacao BBAS4"
#
# Core change: on sheet "1999-2004 (Parte 4)", cell K2's formula changes
# from =J2 to =J2 + '1992-1998 (Parte 5)'!K73, which cascades through the
# whole K/L column (and the IRR/NPV summary cells) via existing formulas.
#
# Also updates the view state: the second sheet ("1999-2004 (Parte 4)")
# becomes the active/selected sheet, with a new selection and scroll
# position on both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1992-1998 (Parte 5)")
$ws2 = $wb.Worksheets.Item("1999-2004 (Parte 4)")

# --- Core data/formula edit -------------------------------------------------
$ws2.Range("K2").Formula = "=J2 + '1992-1998 (Parte 5)'!K73"

# --- Recalculate to cascade dependent formulas ------------------------------
$excel.CalculateFullRebuild()

# --- View state updates ------------------------------------------------------
# Sheet 1: scroll so row 45 is at top, select K2, and it is no longer the
# tab-selected sheet.
$ws1.Activate()
$ws1.Range("K2").Select()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1

# Sheet 2: becomes the active/selected sheet, select K3.
$ws2.Activate()
$ws2.Range("K3").Select()

$wb.Save()
